$d = $word.ActiveDocument

# --- Locate the "Socks in the Dark:" paragraph and the target empty
# --- list paragraph (numId=2) that follows it (skipping the description
# --- paragraph and one blank paragraph).
$rng = $d.Content
$found = $rng.Find.Execute("Socks in the Dark:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $headingPara = $p
        break
    }
}

# Remove the existing _GoBack bookmark that currently sits on the heading
# paragraph ("Socks in the Dark:").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$descriptionPara = $headingPara.Next()
$blankPara = $descriptionPara.Next()
$targetPara = $blankPara.Next()

# --- Question 2: fill in the first empty numbered paragraph.
$targetPara.Range.InsertAfter("This problem you would have to figure out the odds of you picking the right color of sock each time you picked a sock.  The goal would to pull out socks in pairs.  Every time you have to pick 2 socks to keep it consistent. ")

# --- Question 3: new numbered paragraph right after it.
$targetPara.Range.InsertParagraphAfter()
$newPara = $targetPara.Next()
$newPara.Range.InsertAfter("The potential solution could be 4 socks for the first matching pair. The other would be 12 socks to get 3 matching pair. ")

# Re-create the _GoBack bookmark at the very end of the new paragraph's
# text. Zero-length ranges cannot be anchored reliably, so append a
# placeholder character, bookmark it, then delete the placeholder while
# leaving the (now collapsed) bookmark in place.
$newPara.Range.InsertAfter("X")
$endStart = $newPara.Range.End - 2
$placeholder = $d.Range($endStart, $endStart + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder2 = $d.Range($endStart, $endStart + 1)
$placeholder2.Delete()
